$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 29 and 30 swap identity (Coin name + Link columns), then Price/Volume
# columns are updated to their new values for each (post-swap) row.
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  -0.04%  "

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'8.45"
$ws.Range("E30").Value = "'  +2.06%  "

# Remaining rows: refresh Price (D) and Volume(1h) (E) figures.
# A leading apostrophe forces Excel to keep these as literal text (matching
# the existing convention of storing price/volume as strings) instead of
# auto-converting look-alike numbers (e.g. "1.00" -> 1).
$ws.Range("D2").Value = "'63.632.48"
$ws.Range("E2").Value = "'  +1.31%  "
$ws.Range("D3").Value = "'2.656.99"
$ws.Range("E3").Value = "'  +2.85%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'591.46"
$ws.Range("E5").Value = "'  +1.47%  "
$ws.Range("D6").Value = "'144.50"
$ws.Range("E6").Value = "'  -1.09%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E8").Value = "'  -0.78%  "
$ws.Range("D9").Value = "'2.654.37"
$ws.Range("E9").Value = "'  +2.76%  "
$ws.Range("E10").Value = "'  -0.63%  "
$ws.Range("D11").Value = "'5.61"
$ws.Range("E11").Value = "'  -0.21%  "
$ws.Range("E12").Value = "'  +0.46%  "
$ws.Range("E13").Value = "'  +0.34%  "
$ws.Range("D14").Value = "'27.47"
$ws.Range("E14").Value = "'  +1.04%  "
$ws.Range("D15").Value = "'3.129.28"
$ws.Range("E15").Value = "'  +2.70%  "
$ws.Range("D16").Value = "'63.601.30"
$ws.Range("E16").Value = "'  +1.47%  "
$ws.Range("E17").Value = "'  -0.02%  "
$ws.Range("D18").Value = "'2.661.19"
$ws.Range("E18").Value = "'  +3.01%  "
$ws.Range("D19").Value = "'11.45"
$ws.Range("E19").Value = "'  +1.75%  "
$ws.Range("D20").Value = "'340.95"
$ws.Range("E20").Value = "'  -0.06%  "
$ws.Range("E21").Value = "'  -0.28%  "
$ws.Range("D22").Value = "'6.74"
$ws.Range("E22").Value = "'  +0.92%  "
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("D24").Value = "'67.64"
$ws.Range("E24").Value = "'  +0.51%  "
$ws.Range("D25").Value = "'1.65"
$ws.Range("E25").Value = "'  +3.96%  "
$ws.Range("D26").Value = "'1.55"
$ws.Range("E26").Value = "'  +5.88%  "
$ws.Range("E27").Value = "'  -0.17%  "
$ws.Range("D28").Value = "'550.83"
$ws.Range("E28").Value = "'  +18.10%  "
$ws.Range("D31").Value = "'7.80"
$ws.Range("E31").Value = "'  -0.69%  "
$ws.Range("D32").Value = "'1.82"
$ws.Range("E32").Value = "'  +13.63%  "
$ws.Range("E33").Value = "'  +3.31%  "
$ws.Range("D34").Value = "'0.0₃0811"
$ws.Range("E34").Value = "'  -0.71%  "
$ws.Range("D35").Value = "'174.22"
$ws.Range("E35").Value = "'  -1.31%  "
$ws.Range("D36").Value = "'4.91"
$ws.Range("E36").Value = "'  +8.82%  "
$ws.Range("E37").Value = "'  -0.07%  "
$ws.Range("E38").Value = "'  +1.07%  "
$ws.Range("E39").Value = "'  +0.72%  "
$ws.Range("D40").Value = "'1.83"
$ws.Range("E40").Value = "'  +7.51%  "
$ws.Range("D41").Value = "'170.73"
$ws.Range("E41").Value = "'  +7.74%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "'  -0.01%  "
$ws.Range("D43").Value = "'40.28"
$ws.Range("E43").Value = "'  +1.99%  "
$ws.Range("D44").Value = "'3.74"
$ws.Range("E44").Value = "'  +0.07%  "
$ws.Range("D45").Value = "'22.28"
$ws.Range("E45").Value = "'  +6.13%  "
$ws.Range("E46").Value = "'  +0.25%  "
$ws.Range("D47").Value = "'0.0555"
$ws.Range("E47").Value = "'  +2.97%  "
$ws.Range("D48").Value = "'0.0963"
$ws.Range("E48").Value = "'  -0.23%  "
$ws.Range("E49").Value = "'  +1.29%  "
$ws.Range("E50").Value = "'  +2.78%  "
$ws.Range("E51").Value = "'  -0.22%  "
